$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (was "1", now "ლაჩხუთი")
$ws.Name = "ლაჩხუთი"

# Remove the census-note row (old row 2, "(მოსახლეობის აღწერის შედეგებით)"):
# shifts everything below up by one row
$ws.Rows(2).Delete() | Out-Null

# Remove the 1989 and 2002 year columns, keeping only the 2014 column
$ws.Columns("B:C").Delete() | Out-Null

# Match the saved selection state
$ws.Range("A2").Select() | Out-Null
